# "Checked in New tests." — append three new test-case rows (42-44) to the
# "Test Cases" sheet, matching the formatting already used for the rest of
# the table (AuthoringTest!A2 donor style for the bordered data cells, and
# this sheet's own column-B hyperlink-style cell for the JIRA-ID column).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item(1)        # "Test Cases" sheet (sheet1.xml)
$fmt = $wb.Worksheets.Item(3)        # "AuthoringTest" sheet supplies the plain bordered style (s=1)

# --- Copy cell formats first, then overwrite with the real values --------

# Columns A, C, D, E use the plain bordered style already present at
# AuthoringTest!A2 (cellXf index 1 in this workbook).
$fmt.Range("A2").Copy()
$ws.Range("A42:A44").PasteSpecial(-4122)
$ws.Range("C42:E44").PasteSpecial(-4122)

# Column B (JIRA ID) uses the same hyperlink-flavoured style already used
# for every other row in this column.
$ws.Range("B41").Copy()
$ws.Range("B42:B44").PasteSpecial(-4122)

# --- Row 42 : VerifyPostRecordDetails / OPQA-370 --------------------------
$ws.Range("A42").Value = "VerifyPostRecordDetails"
$ws.Range("B42").Value = "OPQA-370"
$ws.Range("C42").Value = "Verify that user contributed articles display the information about the author"
$ws.Range("D42").Value = "Y"
$ws.Range("E42").Value = "PASS"

# --- Row 43 : SeacrhAndViewOwnPost / OPQA-415 ------------------------------
$ws.Range("A43").Value = "SeacrhAndViewOwnPost"
$ws.Range("B43").Value = "OPQA-415"
$ws.Range("C43").Value = "Verify that user is able to search the  posts a user authored themselves and view them."
$ws.Range("D43").Value = "Y"
$ws.Range("E43").Value = "PASS"

# --- Row 44 : SeacrhAndViewOthersPost / OPQA-416 ---------------------------
$ws.Range("A44").Value = "SeacrhAndViewOthersPost"
$ws.Range("B44").Value = "OPQA-416"
$ws.Range("C44").Value = "Verify that user is able to search the posts of others and view them."
$ws.Range("D44").Value = "Y"
$ws.Range("E44").Value = "PASS"

# --- Match the saved view state: last-used cell selected, scrolled down ---
$ws.Range("B44").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow    = 28
$win.ScrollColumn = 1
